$d = $word.ActiveDocument

$pairs = @(
    @("249÷6=41, 3", "226÷4=56, 2"),
    @("771÷6=128, 3", "332÷4=83, 0"),
    @("177÷4=44, 1", "353÷6=58, 5"),
    @("661÷3=220, 1", "877÷8=109, 5"),
    @("451÷5=90, 1", "225÷5=45, 0"),
    @("606÷4=151, 2", "647÷2=323, 1"),
    @("677÷8=84, 5", "230÷2=115, 0"),
    @("442÷7=63, 1", "960÷2=480, 0"),
    @("587÷9=65, 2", "550÷7=78, 4"),
    @("820÷8=102, 4", "894÷2=447, 0"),
    @("645÷5=129, 0", "374÷4=93, 2"),
    @("258÷8=32, 2", "567÷7=81, 0"),
    @("749÷5=149, 4", "759÷7=108, 3"),
    @("344÷7=49, 1", "734÷6=122, 2"),
    @("715÷3=238, 1", "119÷8=14, 7"),
    @("774÷4=193, 2", "436÷8=54, 4"),
    @("677÷6=112, 5", "619÷3=206, 1"),
    @("808÷9=89, 7", "536÷6=89, 2"),
    @("843÷8=105, 3", "874÷2=437, 0"),
    @("745÷9=82, 7", "620÷2=310, 0"),
    @("924÷8=115, 4", "536÷2=268, 0"),
    @("923÷7=131, 6", "502÷6=83, 4"),
    @("377÷3=125, 2", "256÷2=128, 0"),
    @("281÷4=70, 1", "202÷9=22, 4"),
    @("468÷4=117, 0", "573÷6=95, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
